$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values are entered in this precise order so that the resulting
# shared-string table indices (9..21) come out matching the target file
# (Excel appends a new shared string the first time a distinct value is used).
$ws.Range("B6").Value = "more sensitive, fewer false negatives"
$ws.Range("D6").Value = "correctly classifies No  cases more accurately relative to WHO"
$ws.Range("C6").Value = "misses more True Some, both overdiagnoses & underdiagnoses people relative to WHO"
$ws.Range("C8").Value = "sensitivity = correctly says you're not in the second 2 columns (given that you have Severe)"
$ws.Range("C9").Value = "specificity = correctly says you're not in the second 2 rows (given that you have Severe)"
$ws.Range("K9").Value = "adjust parameter " + [char]0x2014 + " assign DALY to undertreating Some"
$ws.Range("K8").Value = "add: one-way sensitivity how bad would penalty have to be for undertreating Some"
$ws.Range("C13").Value = "JP: double check variation in age distribution based on model"
$ws.Range("K10").Value = "state in assumptions that this is not assigned (AL: circle back and do this after rotations, April even just for learning)"
$ws.Range("C14").Value = "emphasize that errors of the models are in the different age groups which affects DALYs"
$ws.Range("C15").Value = "why we used different DALYs for death"
$ws.Range("K13").Value = """Branch H"" (WHO) " + [char]0x2014 + " 56.61046 DALYs conditional on death in Branch H"
$ws.Range("K14").Value = """Branch Q"" (NIRUDAK) " + [char]0x2014 + " 46.11102 DALYs conditional on death in Branch Q"

# --- Formatting: row 6 header cells get wrap text + a taller row ---
$ws.Range("B6:D6").WrapText = $true
$ws.Rows.Item(6).RowHeight = 136

# --- Column C width ---
$ws.Columns.Item(3).ColumnWidth = 13.25

# --- Sheet view: zoom + active cell selection ---
$ws.Application.ActiveWindow.Zoom = 125
[void]$ws.Range("K14").Select()

# --- Reposition the picture, since inserting the rows above pushes it down ---
$shp = $ws.Shapes.Item(1)
$shp.Left = 238.4
$shp.Top = 432.6
$shp.Width = 854.0970078740157
$shp.Height = 312
